$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference range holding the default (unstyled) cell style used
# elsewhere on the sheet, so Price-column text cells keep their look.
$defaultRange = $ws.Range("A1")

# D2: Price column stores numeric-looking values as text
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.001.22'
$ws.Range("D2").Style = $defaultRange.Style
$ws.Range("E2").Value = '  -2.24%  '
# D3: Price column stores numeric-looking values as text
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.831.63'
$ws.Range("D3").Style = $defaultRange.Style
$ws.Range("E3").Value = '  -1.08%  '
$ws.Range("E4").Value = '  -0.02%  '
# D5: Price column stores numeric-looking values as text
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '324.47'
$ws.Range("D5").Style = $defaultRange.Style
$ws.Range("E5").Value = '  -2.80%  '
# D6: Price column stores numeric-looking values as text
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.002'
$ws.Range("D6").Style = $defaultRange.Style
# D7: Price column stores numeric-looking values as text
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4657'
$ws.Range("D7").Style = $defaultRange.Style
$ws.Range("E7").Value = '  -0.12%  '
# D8: Price column stores numeric-looking values as text
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3869'
$ws.Range("D8").Style = $defaultRange.Style
$ws.Range("E8").Value = '  -1.33%  '
# D9: Price column stores numeric-looking values as text
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07880'
$ws.Range("D9").Style = $defaultRange.Style
$ws.Range("E9").Value = '  -0.18%  '
# D10: Price column stores numeric-looking values as text
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.9607'
$ws.Range("D10").Style = $defaultRange.Style
$ws.Range("E10").Value = '  -2.58%  '
# D11: Price column stores numeric-looking values as text
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '21.85'
$ws.Range("D11").Style = $defaultRange.Style
$ws.Range("E11").Value = '  -1.67%  '
# D12: Price column stores numeric-looking values as text
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.856.79'
$ws.Range("D12").Style = $defaultRange.Style
$ws.Range("E12").Value = '  -6.79%  '
# D13: Price column stores numeric-looking values as text
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.670'
$ws.Range("D13").Style = $defaultRange.Style
$ws.Range("E13").Value = '  -3.21%  '
# D14: Price column stores numeric-looking values as text
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.910'
$ws.Range("D14").Style = $defaultRange.Style
$ws.Range("E14").Value = '  -1.65%  '
# D15: Price column stores numeric-looking values as text
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.06828'
$ws.Range("D15").Style = $defaultRange.Style
$ws.Range("E15").Value = '  +0.10%  '
# D16: Price column stores numeric-looking values as text
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '87.32'
$ws.Range("D16").Style = $defaultRange.Style
$ws.Range("E16").Value = '  -0.48%  '
$ws.Range("E17").Value = '  +0.05%  '
# D18: Price column stores numeric-looking values as text
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000009927'
$ws.Range("D18").Style = $defaultRange.Style
$ws.Range("E18").Value = '  -1.53%  '
# D19: Price column stores numeric-looking values as text
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '16.56'
$ws.Range("D19").Style = $defaultRange.Style
$ws.Range("E19").Value = '  -2.76%  '
# D20: Price column stores numeric-looking values as text
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.002'
$ws.Range("D20").Style = $defaultRange.Style
$ws.Range("E20").Value = '  +0.16%  '
# D21: Price column stores numeric-looking values as text
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '28.023.57'
$ws.Range("D21").Style = $defaultRange.Style
$ws.Range("E21").Value = '  -2.22%  '
# D22: Price column stores numeric-looking values as text
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.321'
$ws.Range("D22").Style = $defaultRange.Style
$ws.Range("E22").Value = '  -1.41%  '
# D23: Price column stores numeric-looking values as text
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.97'
$ws.Range("D23").Style = $defaultRange.Style
$ws.Range("E23").Value = '  -2.24%  '
# D24: Price column stores numeric-looking values as text
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.090'
$ws.Range("D24").Style = $defaultRange.Style
$ws.Range("E24").Value = '  -2.00%  '
# D25: Price column stores numeric-looking values as text
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.082.09'
$ws.Range("D25").Style = $defaultRange.Style
$ws.Range("E25").Value = '  -6.24%  '
# D26: Price column stores numeric-looking values as text
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '153.81'
$ws.Range("D26").Style = $defaultRange.Style
$ws.Range("E26").Value = '  +0.12%  '
# D27: Price column stores numeric-looking values as text
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '19.10'
$ws.Range("D27").Style = $defaultRange.Style
$ws.Range("E27").Value = '  -1.63%  '
# D28: Price column stores numeric-looking values as text
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '5.727'
$ws.Range("D28").Style = $defaultRange.Style
$ws.Range("E28").Value = '  -5.90%  '
# D29: Price column stores numeric-looking values as text
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.966'
$ws.Range("D29").Style = $defaultRange.Style
$ws.Range("E29").Value = '  -2.75%  '
# D30: Price column stores numeric-looking values as text
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '117.60'
$ws.Range("D30").Style = $defaultRange.Style
$ws.Range("E30").Value = '  +0.11%  '
# D31: Price column stores numeric-looking values as text
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.09264'
$ws.Range("D31").Style = $defaultRange.Style
$ws.Range("E31").Value = '  -1.60%  '
# D32: Price column stores numeric-looking values as text
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.9338'
$ws.Range("D32").Style = $defaultRange.Style
$ws.Range("E32").Value = '  -4.71%  '
# D33: Price column stores numeric-looking values as text
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.284'
$ws.Range("D33").Style = $defaultRange.Style
$ws.Range("E33").Value = '  -1.57%  '
# D34: Price column stores numeric-looking values as text
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.317'
$ws.Range("D34").Style = $defaultRange.Style
$ws.Range("E34").Value = '  -2.23%  '
# D35: Price column stores numeric-looking values as text
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.293'
$ws.Range("D35").Style = $defaultRange.Style
$ws.Range("E35").Value = '  -5.86%  '
# D36: Price column stores numeric-looking values as text
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.05868'
$ws.Range("D36").Style = $defaultRange.Style
$ws.Range("E36").Value = '  -4.12%  '
# D37: Price column stores numeric-looking values as text
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02143'
$ws.Range("D37").Style = $defaultRange.Style
$ws.Range("E37").Value = '  -2.57%  '
# D38: Price column stores numeric-looking values as text
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.143'
$ws.Range("D38").Style = $defaultRange.Style
$ws.Range("E38").Value = '  -2.00%  '
# D39: Price column stores numeric-looking values as text
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '7.779'
$ws.Range("D39").Style = $defaultRange.Style
$ws.Range("E39").Value = '  +2.36%  '
# D40: Price column stores numeric-looking values as text
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.5580'
$ws.Range("D40").Style = $defaultRange.Style
$ws.Range("E40").Value = '  -2.22%  '
# D41: Price column stores numeric-looking values as text
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '9.864'
$ws.Range("D41").Style = $defaultRange.Style
$ws.Range("E41").Value = '  -2.29%  '
# D42: Price column stores numeric-looking values as text
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1759'
$ws.Range("D42").Style = $defaultRange.Style
$ws.Range("E42").Value = '  -1.82%  '
# D43: Price column stores numeric-looking values as text
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '11.58'
$ws.Range("D43").Style = $defaultRange.Style
$ws.Range("E43").Value = '  -1.99%  '
# D44: Price column stores numeric-looking values as text
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.5253'
$ws.Range("D44").Style = $defaultRange.Style
$ws.Range("E44").Value = '  -2.52%  '
# D45: Price column stores numeric-looking values as text
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.07012'
$ws.Range("D45").Style = $defaultRange.Style
$ws.Range("E45").Value = '  -2.08%  '
# D46: Price column stores numeric-looking values as text
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.138'
$ws.Range("D46").Style = $defaultRange.Style
$ws.Range("E46").Value = '  -10.95%  '
$ws.Range("E47").Value = '  -4.32%  '
# D48: Price column stores numeric-looking values as text
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '113.26'
$ws.Range("D48").Style = $defaultRange.Style
$ws.Range("E48").Value = '  -0.08%  '
# D49: Price column stores numeric-looking values as text
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.102'
$ws.Range("D49").Style = $defaultRange.Style
$ws.Range("E49").Value = '  -12.10%  '
$ws.Range("E50").Value = '  +0.02%  '
# D51: Price column stores numeric-looking values as text
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.321'
$ws.Range("D51").Style = $defaultRange.Style
$ws.Range("E51").Value = '  +0.61%  '
